$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set Obrigatorio (column E) to "S" for rows 2-9 and 11-15 (row 10 stays "N")
$rows = @(2,3,4,5,6,7,8,9,11,12,13,14,15)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 5).Value = "S"
}
